# Generate Report for Handback
# The "db9b99e9-93b7-4156-a7fd-3b35655629bf.md" file has finished its
# handback round-trip: flip its status from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is reported, stamp the
# new handback timestamps per-locale, and clear the stale "not latest"
# error detail now that the handback is current.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet: one row per source file, one status column per locale ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusDone   # zh-cn status
$wsOverview.Range("F3").Value = $statusDone   # de-de status

# --- zh-cn sheet: detail row for db9b99e9...md ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusDone                     # Status
$wsZh.Range("K3").Value = "2016-08-20 10:53:28"           # Latest Handback DateTime
$wsZh.Range("P3").Value = ""                              # Error Detail cleared
$wsZh.Columns.Item(16).ColumnWidth = 12.75                # Error Detail column shrinks back down

# --- de-de sheet: detail row for db9b99e9...md ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusDone                     # Status
$wsDe.Range("K3").Value = "2016-08-20 10:53:33"           # Latest Handback DateTime
$wsDe.Range("P3").Value = ""                              # Error Detail cleared
$wsDe.Columns.Item(16).ColumnWidth = 12.75                # Error Detail column shrinks back down
